$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Saturday date strings (rows 2-12) in columns A and B
$satOld = "Saturday Saturday_ February 4_ 2023"
$satNew = "Saturday Saturday_ February 2023"
for ($r = 2; $r -le 12; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    if ($cellA.Value() -eq $satOld) {
        $cellA.Value = $satNew
    }
    $cellB = $ws.Cells.Item($r, 2)
    if ($cellB.Value() -eq $satOld) {
        $cellB.Value = $satNew
    }
}

# Update Sunday date strings (rows 13-23) in columns A and B
$sunOld = "Sunday Sunday_ February 6_ 2023"
$sunNew = "Sunday Sunday_ February 2023"
for ($r = 13; $r -le 23; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    if ($cellA.Value() -eq $sunOld) {
        $cellA.Value = $sunNew
    }
    $cellB = $ws.Cells.Item($r, 2)
    if ($cellB.Value() -eq $sunOld) {
        $cellB.Value = $sunNew
    }
}

# Clear the erroneous ERROR value in column C for rows 3, 12, 13
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(12, 3).Value = ""
$ws.Cells.Item(13, 3).Value = ""

# Fix "ISAAC S SANCHEZ" -> "ISAAC SANCHEZ" in rows 7 and 18, column C
$ws.Cells.Item(7, 3).Value = "ISAAC SANCHEZ"
$ws.Cells.Item(18, 3).Value = "ISAAC SANCHEZ"
